# Daily_Log.xlsx update:
#  - June 21-2023: record the new accomplishment in column B next to the
#    "Committed By:" / name rows (row 12), and move the sheet's selection.
#  - Append four new daily-log sheets (June 24-2023 .. June 27-2023),
#    each a copy of the still-blank "June 22-2023" template sheet.

$wb = $excel.ActiveWorkbook

# --- Update "June 21-2023": add the new log entry -------------------------
$ws21 = $wb.Worksheets.Item("June 21-2023")
$ws21.Range("B12").Value = "added functionality of food exchange computation."

# --- Append the four new blank daily-log sheets ----------------------------
# "June 22-2023" is the blank template (column A filled in, column B empty)
# that the new days are cloned from.
$template = $wb.Worksheets.Item("June 22-2023")
$newNames = @("June 24-2023", "June 25-2023", "June 26-2023", "June 27-2023")

foreach ($name in $newNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $name
    $newSheet.Range("D25").Select()
}

# Last new sheet (June 27-2023) ends up with a different selected cell.
$wb.Worksheets.Item("June 27-2023").Range("L26").Select()

# Restore the originally active sheet/selection on "June 21-2023".
$ws21.Activate()
$ws21.Range("B17").Select()
